$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.365.76'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.786.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.22%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.552'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.80%  '

$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.73'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.57%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.281'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0660'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.83%  '

$ws.Range("E11").Value = '  +0.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.042.68'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.21%  '

$ws.Range("E13").Value = '  +6.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.770.16'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.634'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.97%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.374.15'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.03'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '254.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.28%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0745'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.86%  '

$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.91%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.31%  '

$ws.Range("E24").Value = '  -4.54%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.55'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.35%  '

$ws.Range("E28").Value = '  -3.86%  '

$ws.Range("E29").Value = '  +0.25%  '

$ws.Range("E30").Value = '  -3.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0513'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.17%  '

$ws.Range("E32").Value = '  -2.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.60'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.35%  '

$ws.Range("E34").Value = '  +3.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.452.02'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.630'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.35%  '

$ws.Range("E38").Value = '  -1.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '83.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.85'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.21%  '

$ws.Range("E41").Value = '  -0.32%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.889'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.99%  '

$ws.Range("E43").Value = '  -2.50%  '

$ws.Range("E44").Value = '  -4.16%  '

$ws.Range("E45").Value = '  -2.32%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.86'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.942.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.24'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.27%  '

$ws.Range("E49").Value = '  +0.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '99.07'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.00%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '50.13'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.31%  '
